$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Content change -------------------------------------------------
# A new "jan" column header is added to the example2 sheet (column J,
# row 5), next to the existing "Timeseries b in euros" header in C5.
# This is the substantive edit behind the commit: read_ts_xlsx now
# reads exactly as many columns as the determined column types, and
# this fixture gained one more recognised column ("jan").
$ws2.Range("J5").Value = "jan"

# --- View / selection state ------------------------------------------
# Keep "example2" the active sheet (matches original activeTab/tabSelected)
# and move its selection the way the source workbook shows it.
$ws2.Activate()
$ws2.Range("I25").Select()

# --- Cosmetic column-width touch-ups ----------------------------------
# The canonical file shows small width reductions across both sheets'
# column metadata (a side effect of the resave that produced the new
# column). Reproduce the per-column widths as closely as the host
# allows.
$ws1.Range($ws1.Cells.Item(1, 1), $ws1.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 11.044217687074866

$ws2.Columns.Item(1).ColumnWidth = 23.467687074829968
$ws2.Columns.Item(2).ColumnWidth = 1.5952380952380967
$ws2.Columns.Item(3).ColumnWidth = 7.671768707482998
$ws2.Columns.Item(4).ColumnWidth = 33.722789115646265
$ws2.Columns.Item(5).ColumnWidth = 8.212585034013607
$ws2.Columns.Item(6).ColumnWidth = 12.528911564625867
$ws2.Range($ws2.Cells.Item(1, 7), $ws2.Cells.Item(1, 8)).EntireColumn.ColumnWidth = 8.212585034013607
$ws2.Range($ws2.Cells.Item(1, 9), $ws2.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 12.528911564625867
